$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "wins_data"
$ws.Cells.Item(2,3).Value = 88
$ws.Cells.Item(2,4).Value = 5
$ws.Cells.Item(2,5).Value = 3126.129388809204
$ws.Cells.Item(2,6).Value = 5
$ws.Cells.Item(2,7).Value = 15

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "wins_data"
$ws.Cells.Item(3,3).Value = 87
$ws.Cells.Item(3,4).Value = 7
$ws.Cells.Item(3,5).Value = 2639.851808547974
$ws.Cells.Item(3,6).Value = 2
$ws.Cells.Item(3,7).Value = 13

$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "wins_data"
$ws.Cells.Item(4,3).Value = 90
$ws.Cells.Item(4,4).Value = 10
$ws.Cells.Item(4,5).Value = 2681.691646575928
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 10

$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "wins_data"
$ws.Cells.Item(5,3).Value = 90
$ws.Cells.Item(5,4).Value = 8
$ws.Cells.Item(5,5).Value = 2624.805212020874
$ws.Cells.Item(5,6).Value = 3
$ws.Cells.Item(5,7).Value = 12

$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "wins_data"
$ws.Cells.Item(6,3).Value = 89
$ws.Cells.Item(6,4).Value = 13
$ws.Cells.Item(6,5).Value = 2650.022268295288
$ws.Cells.Item(6,6).Value = 2
$ws.Cells.Item(6,7).Value = 7

$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "wins_data"
$ws.Cells.Item(7,3).Value = 90
$ws.Cells.Item(7,4).Value = 9
$ws.Cells.Item(7,5).Value = 2633.105278015137
$ws.Cells.Item(7,6).Value = 3
$ws.Cells.Item(7,7).Value = 11

$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "wins_data"
$ws.Cells.Item(8,3).Value = 92
$ws.Cells.Item(8,4).Value = 10
$ws.Cells.Item(8,5).Value = 2665.648937225342
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 10

$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "wins_data"
$ws.Cells.Item(9,3).Value = 91
$ws.Cells.Item(9,4).Value = 12
$ws.Cells.Item(9,5).Value = 2632.755041122437
$ws.Cells.Item(9,6).Value = 0
$ws.Cells.Item(9,7).Value = 8

$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "wins_data"
$ws.Cells.Item(10,3).Value = 90
$ws.Cells.Item(10,4).Value = 12
$ws.Cells.Item(10,5).Value = 2627.527475357056
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 8

$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "wins_data"
$ws.Cells.Item(11,3).Value = 91
$ws.Cells.Item(11,4).Value = 9
$ws.Cells.Item(11,5).Value = 2648.579359054565
$ws.Cells.Item(11,6).Value = 0
$ws.Cells.Item(11,7).Value = 11

$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "wins_data"
$ws.Cells.Item(12,3).Value = 91
$ws.Cells.Item(12,4).Value = 12
$ws.Cells.Item(12,5).Value = 2614.293813705444
$ws.Cells.Item(12,6).Value = 0
$ws.Cells.Item(12,7).Value = 8

$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "wins_data"
$ws.Cells.Item(13,3).Value = 90
$ws.Cells.Item(13,4).Value = 13
$ws.Cells.Item(13,5).Value = 2633.073806762695
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 7

$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "wins_data"
$ws.Cells.Item(14,3).Value = 93
$ws.Cells.Item(14,4).Value = 12
$ws.Cells.Item(14,5).Value = 2650.115728378296
$ws.Cells.Item(14,6).Value = 0
$ws.Cells.Item(14,7).Value = 8

$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "wins_data"
$ws.Cells.Item(15,3).Value = 91
$ws.Cells.Item(15,4).Value = 14
$ws.Cells.Item(15,5).Value = 2613.082647323608
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 6

$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "wins_data"
$ws.Cells.Item(16,3).Value = 91
$ws.Cells.Item(16,4).Value = 15
$ws.Cells.Item(16,5).Value = 2673.084735870361
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 5

$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "wins_data"
$ws.Cells.Item(17,3).Value = 92
$ws.Cells.Item(17,4).Value = 12
$ws.Cells.Item(17,5).Value = 2675.118923187256
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 8

$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "wins_data"
$ws.Cells.Item(18,3).Value = 93
$ws.Cells.Item(18,4).Value = 13
$ws.Cells.Item(18,5).Value = 2635.033130645752
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = 7

$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = "wins_data"
$ws.Cells.Item(19,3).Value = 93
$ws.Cells.Item(19,4).Value = 12
$ws.Cells.Item(19,5).Value = 2617.598533630371
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = 8

$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = "wins_data"
$ws.Cells.Item(20,3).Value = 90
$ws.Cells.Item(20,4).Value = 13
$ws.Cells.Item(20,5).Value = 2600.1136302948
$ws.Cells.Item(20,6).Value = 1
$ws.Cells.Item(20,7).Value = 7

$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = "wins_data"
$ws.Cells.Item(21,3).Value = 93
$ws.Cells.Item(21,4).Value = 13
$ws.Cells.Item(21,5).Value = 2647.555112838745
$ws.Cells.Item(21,6).Value = 0
$ws.Cells.Item(21,7).Value = 7

$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = "wins_data"
$ws.Cells.Item(22,3).Value = 93
$ws.Cells.Item(22,4).Value = 13
$ws.Cells.Item(22,5).Value = 2670.164346694946
$ws.Cells.Item(22,6).Value = 1
$ws.Cells.Item(22,7).Value = 7

$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = "wins_data"
$ws.Cells.Item(23,3).Value = 93
$ws.Cells.Item(23,4).Value = 13
$ws.Cells.Item(23,5).Value = 2620.767116546631
$ws.Cells.Item(23,6).Value = 1
$ws.Cells.Item(23,7).Value = 7
